# Auto-generated edit script: refresh market-price-derived columns
# (currentAveragePrice / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ)
# across the per-crafting-class sheets, per the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 234872.55
$ws.Range("I106").Value = 329484.94
$ws.Range("K106").Value = 329484.94
$ws.Range("M106").Value = -328853.94
$ws.Range("H111").Value = 2425
$ws.Range("I111").Value = 1666.6666
$ws.Range("J111").Value = 2880
$ws.Range("K111").Value = 4999.9998
$ws.Range("L111").Value = 8640
$ws.Range("M111").Value = -1932.9998
$ws.Range("N111").Value = -14774
$ws.Range("H115").Value = 679.64703
$ws.Range("I115").Value = 709.625
$ws.Range("J115").Value = 200
$ws.Range("K115").Value = 2128.875
$ws.Range("L115").Value = 600
$ws.Range("M115").Value = -561.875
$ws.Range("N115").Value = -3734
$ws.Range("H135").Value = 650.1429000000001
$ws.Range("I135").Value = 333.57144
$ws.Range("J135").Value = 1599.8572
$ws.Range("K135").Value = 3002.14296
$ws.Range("L135").Value = 14398.7148
$ws.Range("M135").Value = -467.1429600000001
$ws.Range("N135").Value = -19468.7148

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1157.2
$ws.Range("I2").Value = 539.8570999999999
$ws.Range("K2").Value = 539.8570999999999
$ws.Range("M2").Value = -426.8570999999999
$ws.Range("H45").Value = 1612.4615
$ws.Range("I45").Value = 1358.909
$ws.Range("K45").Value = 1358.909
$ws.Range("M45").Value = -981.9090000000001
$ws.Range("H74").Value = 2993
$ws.Range("I74").Value = 1757.3158
$ws.Range("K74").Value = 1757.3158
$ws.Range("M74").Value = -883.3158000000001
$ws.Range("H77").Value = 2993
$ws.Range("I77").Value = 1757.3158
$ws.Range("K77").Value = 8786.579
$ws.Range("M77").Value = -4418.579
$ws.Range("H110").Value = 31117.947
$ws.Range("I110").Value = 44849.383
$ws.Range("K110").Value = 44849.383
$ws.Range("M110").Value = -42804.383
$ws.Range("H116").Value = 1157.2
$ws.Range("I116").Value = 539.8570999999999
$ws.Range("K116").Value = 539.8570999999999
$ws.Range("M116").Value = 1754.1429
$ws.Range("H122").Value = 3878.45
$ws.Range("I122").Value = 3959.2693
$ws.Range("J122").Value = 3728.3572
$ws.Range("K122").Value = 11877.8079
$ws.Range("L122").Value = 11185.0716
$ws.Range("M122").Value = -9427.8079
$ws.Range("N122").Value = -16085.0716
$ws.Range("H132").Value = 2598.0256
$ws.Range("I132").Value = 2104.7407
$ws.Range("J132").Value = 3707.9167
$ws.Range("K132").Value = 6314.222099999999
$ws.Range("L132").Value = 11123.7501
$ws.Range("M132").Value = -3784.222099999999
$ws.Range("N132").Value = -16183.7501
$ws.Range("H134").Value = 90000
$ws.Range("J134").Value = 90000
$ws.Range("L134").Value = 90000
$ws.Range("N134").Value = -100140
$ws.Range("H140").Value = 155520.5
$ws.Range("J140").Value = 155520.5
$ws.Range("L140").Value = 155520.5
$ws.Range("N140").Value = -165880.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1157.2
$ws.Range("I3").Value = 539.8570999999999
$ws.Range("K3").Value = 539.8570999999999
$ws.Range("M3").Value = -425.8570999999999
$ws.Range("H86").Value = 1967.7368
$ws.Range("I86").Value = 1529.8462
$ws.Range("J86").Value = 2916.5
$ws.Range("K86").Value = 1529.8462
$ws.Range("L86").Value = 2916.5
$ws.Range("M86").Value = -406.8462
$ws.Range("N86").Value = -5162.5
$ws.Range("H89").Value = 1967.7368
$ws.Range("I89").Value = 1529.8462
$ws.Range("J89").Value = 2916.5
$ws.Range("K89").Value = 7649.231
$ws.Range("L89").Value = 14582.5
$ws.Range("M89").Value = -2033.231
$ws.Range("N89").Value = -25814.5
$ws.Range("H105").Value = 1786.6666
$ws.Range("I105").Value = 1493.5
$ws.Range("J105").Value = 2373
$ws.Range("K105").Value = 1493.5
$ws.Range("L105").Value = 2373
$ws.Range("M105").Value = 253.5
$ws.Range("N105").Value = -5867

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 71432050
$ws.Range("I86").Value = 111114940
$ws.Range("K86").Value = 111114940
$ws.Range("M86").Value = -111113817
$ws.Range("H89").Value = 71432050
$ws.Range("I89").Value = 111114940
$ws.Range("K89").Value = 555574700
$ws.Range("M89").Value = -555569084
$ws.Range("H135").Value = 65332.11
$ws.Range("J135").Value = 65332.11
$ws.Range("L135").Value = 65332.11
$ws.Range("N135").Value = -75472.11

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("H107").Value = 437.25
$ws.Range("I107").Value = 833.3333
$ws.Range("J107").Value = 345.84616
$ws.Range("K107").Value = 2499.9999
$ws.Range("L107").Value = 1037.53848
$ws.Range("M107").Value = -579.9998999999998
$ws.Range("N107").Value = -4877.53848
$ws.Range("H116").Value = 2922.1428
$ws.Range("I116").Value = 1511
$ws.Range("K116").Value = 4533
$ws.Range("M116").Value = -1091
$ws.Range("H132").Value = 4082.9666
$ws.Range("I132").Value = 4071.5
$ws.Range("J132").Value = 4084.7307
$ws.Range("K132").Value = 36643.5
$ws.Range("L132").Value = 36762.5763
$ws.Range("M132").Value = -34113.5
$ws.Range("N132").Value = -41822.5763

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 16383
$ws.Range("J96").Value = 16612.2
$ws.Range("L96").Value = 16612.2
$ws.Range("N96").Value = -22104.2
$ws.Range("H102").Value = 963.2759
$ws.Range("I102").Value = 1005.7692
$ws.Range("J102").Value = 595
$ws.Range("K102").Value = 1005.7692
$ws.Range("L102").Value = 595
$ws.Range("M102").Value = 616.2308
$ws.Range("N102").Value = -3839

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1276.0333
$ws.Range("I16").Value = 1299.3334
$ws.Range("K16").Value = 1299.3334
$ws.Range("M16").Value = -1129.3334
$ws.Range("H61").Value = 672.75
$ws.Range("I61").Value = 672.75
$ws.Range("K61").Value = 672.75
$ws.Range("M61").Value = -470.75
$ws.Range("H68").Value = 3570.8572
$ws.Range("I68").Value = 3499.5
$ws.Range("J68").Value = 3999
$ws.Range("K68").Value = 3499.5
$ws.Range("L68").Value = 3999
$ws.Range("M68").Value = -2750.5
$ws.Range("N68").Value = -5497
$ws.Range("H71").Value = 3570.8572
$ws.Range("I71").Value = 3499.5
$ws.Range("J71").Value = 3999
$ws.Range("K71").Value = 17497.5
$ws.Range("L71").Value = 19995
$ws.Range("M71").Value = -13753.5
$ws.Range("N71").Value = -27483
$ws.Range("H113").Value = 672.75
$ws.Range("I113").Value = 672.75
$ws.Range("K113").Value = 672.75
$ws.Range("M113").Value = 1497.25
$ws.Range("H134").Value = 49000
$ws.Range("J134").Value = 50000
$ws.Range("L134").Value = 50000
$ws.Range("N134").Value = -60140
$ws.Range("H141").Value = 112993.8
$ws.Range("J141").Value = 112993.8
$ws.Range("L141").Value = 112993.8
$ws.Range("N141").Value = -123353.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 143299.11
$ws.Range("J46").Value = 143299.11
$ws.Range("L46").Value = 143299.11
$ws.Range("N46").Value = -143761.11
$ws.Range("H122").Value = 4563.909
$ws.Range("I122").Value = 3031.1
$ws.Range("J122").Value = 5841.25
$ws.Range("K122").Value = 9093.299999999999
$ws.Range("L122").Value = 17523.75
$ws.Range("M122").Value = -6643.299999999999
$ws.Range("N122").Value = -22423.75
$ws.Range("H133").Value = 79999.336
$ws.Range("J133").Value = 79999.336
$ws.Range("L133").Value = 79999.336
$ws.Range("N133").Value = -90119.336
$ws.Range("H134").Value = 143299.11
$ws.Range("J134").Value = 143299.11
$ws.Range("L134").Value = 429897.33
$ws.Range("N134").Value = -434967.33
$ws.Range("H138").Value = 108993.664
$ws.Range("J138").Value = 108993.664
$ws.Range("L138").Value = 108993.664
$ws.Range("N138").Value = -119273.664

